$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stock-card entry row (25.08.2022): opening balance carried from
# previous balance (F12), an issue of 4, and the resulting balance.
$ws.Range("A13").Value = "25.08.2022"
$ws.Range("D13").Value = 4
$ws.Range("F13").Formula = "=B13-D13"
$ws.Range("B14").Formula = "=F13"

# Row 11 no longer needs an explicit custom height.
$ws.Rows.Item(11).AutoFit()

# Update the active selection shown when the workbook is reopened.
[void]$ws.Range("B13:C14").Select()
